$wb = $excel.ActiveWorkbook

# Rename the existing (only) sheet to "Stock Finances"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Stock Finances"

# Add a second sheet named "Current Price Data", placed after the first sheet
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Current Price Data"

# Fill in the "Stock Finances" sheet data
$data = @(
    @("Stock Name", "Symbol", "Earnings Estimate", "Revenue Estimates"),
    @("Apple", "Apple Inc. (AAPL)", "Avg. Estimate", "Avg. Estimate"),
    @("Tesla", "Tesla, Inc. (TSLA)", "Avg. Estimate", "Avg. Estimate"),
    @("Microsoft", "Microsoft Corporation (MSFT)", "Avg. Estimate", "Avg. Estimate"),
    @("Sony", "Sony Group Corporation (SONY)", "Avg. Estimate", "Avg. Estimate"),
    @("GameStop", "GameStop Corp. (GME)", "Avg. Estimate", "Avg. Estimate"),
    @("Virgin Galactic", "Virgin Galactic Holdings, Inc. (SPCE)", "Avg. Estimate", "Avg. Estimate"),
    @("Honda", "Honda Motor Co., Ltd. (HMC)", "Avg. Estimate", "Avg. Estimate"),
    @("Toyota", "Toyota Motor Corporation (TM)", "Avg. Estimate", "Avg. Estimate"),
    @("Boeing", "The Boeing Company (BA)", "Avg. Estimate", "Avg. Estimate"),
    @("Nintendo", "Nintendo Co., Ltd. (NTDOY)", "Avg. Estimate", "Avg. Estimate"),
    @("T-Mobile", "T-Mobile US, Inc. (TMUS)", "Avg. Estimate", "Avg. Estimate")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws1.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

# Auto-fit columns A:D so widths match the bestFit columns in the target
$ws1.Range("A1:D12").EntireColumn.AutoFit() | Out-Null

# Re-apply default formatting attributes across the populated range; this mirrors
# the author's paste-derived cell style (applyNumberFormat/applyFill/applyAlignment/
# applyProtection all flagged, still resolving to the default look) without changing
# the visible appearance of the data.
$dataRange = $ws1.Range("A1:D12")
$dataRange.WrapText = $false

# Select B2 on the "Stock Finances" sheet and make it the active/tab-selected sheet
$ws1.Range("B2").Select()
$ws1.Activate()
